# Config.py test data sheet: add TC02 test steps (rows 7-12), widen column B,
# add a new "vertical top + wrap" cell style, and update the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width (16.140625 -> 20.7109375 chars) ---------------------
$ws.Columns.Item(2).ColumnWidth = 19.9

# --- Row 7: click_admin ---------------------------------------------------
$ws.Range("A7").Value = "TC02"
$ws.Range("A7").WrapText = $true
$ws.Range("A7").VerticalAlignment = -4108

$ws.Range("B7").Value = "click_admin"
$ws.Range("B7").VerticalAlignment = -4108

$ws.Range("C7").Value = "xpath"
$ws.Range("C7").WrapText = $true
$ws.Range("C7").VerticalAlignment = -4108

$ws.Range("D7").Value = "//span[@class='oxd-text oxd-text--span oxd-main-menu-item--name'][normalize-space()='Admin']"
$ws.Range("D7").WrapText = $true

$ws.Rows.Item(7).RowHeight = 45

# --- Row 8: click_job -------------------------------------------------
$ws.Range("A8").Value = "TC02"
$ws.Range("A8").WrapText = $true
$ws.Range("A8").VerticalAlignment = -4108

$ws.Range("B8").Value = "click_job"

$ws.Range("C8").Value = "xpath"
$ws.Range("C8").WrapText = $true
$ws.Range("C8").VerticalAlignment = -4108

$ws.Range("D8").Value = "//span[normalize-space()='Job']"

$ws.Rows.Item(8).RowHeight = 30

# --- Row 9: click_employmentstatus ----------------------------------------
$ws.Range("A9").Value = "TC02"
$ws.Range("A9").WrapText = $true
$ws.Range("A9").VerticalAlignment = -4108

$ws.Range("B9").Value = "click_employmentstatus"
$ws.Range("B9").VerticalAlignment = -4108

$ws.Range("C9").Value = "xpath"
$ws.Range("C9").WrapText = $true
$ws.Range("C9").VerticalAlignment = -4108

$ws.Range("D9").Value = "//a[text()='Employment Status']"

# --- Row 10: click_add ------------------------------------------------
$ws.Range("A10").Value = "TC02"
$ws.Range("A10").WrapText = $true
$ws.Range("A10").VerticalAlignment = -4108

$ws.Range("B10").Value = "click_add"
$ws.Range("B10").VerticalAlignment = -4108

$ws.Range("C10").Value = "xpath"
$ws.Range("C10").WrapText = $true
$ws.Range("C10").VerticalAlignment = -4108

$ws.Range("D10").Value = "//button[normalize-space()='Add']"
$ws.Range("D10").WrapText = $true
$ws.Range("D10").HorizontalAlignment = -4131

$ws.Rows.Item(10).RowHeight = 60

# --- Row 11: enter_name -------------------------------------------------
$ws.Range("A11").Value = "TC02"
$ws.Range("A11").WrapText = $true
$ws.Range("A11").VerticalAlignment = -4108

$ws.Range("B11").Value = "enter_name"
$ws.Range("B11").VerticalAlignment = -4108

$ws.Range("C11").Value = "xpath"
$ws.Range("C11").WrapText = $true
$ws.Range("C11").VerticalAlignment = -4108

$ws.Range("D11").Value = "//div[@class='oxd-input-group oxd-input-field-bottom-space']//div//input[@class='oxd-input oxd-input--active']"
$ws.Range("D11").WrapText = $true

$ws.Range("E11").Value = "Testing"

$ws.Rows.Item(11).RowHeight = 30

# --- Row 12: click_save -------------------------------------------------
$ws.Range("A12").Value = "TC02"
$ws.Range("A12").WrapText = $true
$ws.Range("A12").VerticalAlignment = -4108

$ws.Range("B12").Value = "click_save"
$ws.Range("B12").VerticalAlignment = -4108

$ws.Range("C12").Value = "xpath"
$ws.Range("C12").WrapText = $true
$ws.Range("C12").VerticalAlignment = -4108

$ws.Range("D12").Value = "//button[normalize-space()='Save']"
$ws.Range("D12").WrapText = $true
$ws.Range("D12").VerticalAlignment = -4160

$ws.Rows.Item(12).RowHeight = 75

# --- Sheet view: scroll so row 8 is at the top, select D15 --------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D15").Select() | Out-Null
